$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 134 <- data from former row 137
$ws.Range("B134").Value = 7483188
$ws.Range("C134").Value = "Ecuador LigaPro Serie A"
$ws.Range("D134").Value = 45256.83333333334
$ws.Range("E134").Value = "Gualaceo SC"
$ws.Range("F134").Value = "Emelec"
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 2
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 1
$ws.Range("K134").Value = "A"
$ws.Range("L134").Value = 3.6
$ws.Range("M134").Value = 3.3
$ws.Range("N134").Value = 2.05
$ws.Range("O134").Value = 2.6
$ws.Range("P134").Value = 3.25
$ws.Range("Q134").Value = 2.75
$ws.Range("R134").Value = 0
$ws.Range("S134").Value = 1.8
$ws.Range("T134").Value = 2
$ws.Range("U134").Value = 2.5
$ws.Range("V134").Value = 1.975
$ws.Range("W134").Value = 1.825
$ws.Range("X134").Value = -1
$ws.Range("Y134").Value = -1
$ws.Range("Z134").Value = 1.75
$ws.Range("AA134").Value = -1
$ws.Range("AB134").Value = 1
$ws.Range("AC134").Value = -1
$ws.Range("AD134").Value = 0.825

# Row 135 <- data from former row 136
$ws.Range("B135").Value = 7482867
$ws.Range("C135").Value = "Ecuador LigaPro Serie A"
$ws.Range("D135").Value = 45256.83333333334
$ws.Range("E135").Value = "Cumbaya FC"
$ws.Range("F135").Value = "LDU Quito"
$ws.Range("G135").Value = 1
$ws.Range("H135").Value = 2
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = "A"
$ws.Range("L135").Value = 5.25
$ws.Range("M135").Value = 3.75
$ws.Range("N135").Value = 1.65
$ws.Range("O135").Value = 9
$ws.Range("P135").Value = 4.5
$ws.Range("Q135").Value = 1.363
$ws.Range("R135").Value = 1.25
$ws.Range("S135").Value = 1.975
$ws.Range("T135").Value = 1.825
$ws.Range("U135").Value = 2.5
$ws.Range("V135").Value = 1.825
$ws.Range("W135").Value = 1.975
$ws.Range("X135").Value = -1
$ws.Range("Y135").Value = -1
$ws.Range("Z135").Value = 0.363
$ws.Range("AA135").Value = 0.4875
$ws.Range("AB135").Value = -0.5
$ws.Range("AC135").Value = 0.825
$ws.Range("AD135").Value = -1

# Row 136 <- data from former row 135
$ws.Range("B136").Value = 7482832
$ws.Range("C136").Value = "Ecuador LigaPro Serie A"
$ws.Range("D136").Value = 45256.83333333334
$ws.Range("E136").Value = "Barcelona Guayaquil"
$ws.Range("F136").Value = "Guayaquil City"
$ws.Range("G136").Value = 2
$ws.Range("H136").Value = 1
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = "H"
$ws.Range("L136").Value = 1.363
$ws.Range("M136").Value = 5
$ws.Range("N136").Value = 7.5
$ws.Range("O136").Value = 1.444
$ws.Range("P136").Value = 4
$ws.Range("Q136").Value = 8
$ws.Range("R136").Value = -1.25
$ws.Range("S136").Value = 2.05
$ws.Range("T136").Value = 1.75
$ws.Range("U136").Value = 2.5
$ws.Range("V136").Value = 1.95
$ws.Range("W136").Value = 1.85
$ws.Range("X136").Value = 0.444
$ws.Range("Y136").Value = -1
$ws.Range("Z136").Value = -1
$ws.Range("AA136").Value = -0.5
$ws.Range("AB136").Value = 0.375
$ws.Range("AC136").Value = 0.95
$ws.Range("AD136").Value = -1

# Row 137 <- data from former row 134
$ws.Range("B137").Value = 7483306
$ws.Range("C137").Value = "Ecuador LigaPro Serie A"
$ws.Range("D137").Value = 45256.83333333334
$ws.Range("E137").Value = "Tecnico Universitario"
$ws.Range("F137").Value = "Club Atletico Libertad"
$ws.Range("G137").Value = 1
$ws.Range("H137").Value = 1
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = "D"
$ws.Range("L137").Value = 1.5
$ws.Range("M137").Value = 4.333
$ws.Range("N137").Value = 5.75
$ws.Range("O137").Value = 1.533
$ws.Range("P137").Value = 4.2
$ws.Range("Q137").Value = 5.5
$ws.Range("R137").Value = -1
$ws.Range("S137").Value = 1.925
$ws.Range("T137").Value = 1.875
$ws.Range("U137").Value = 2.25
$ws.Range("V137").Value = 1.8
$ws.Range("W137").Value = 2
$ws.Range("X137").Value = -1
$ws.Range("Y137").Value = 3.2
$ws.Range("Z137").Value = -1
$ws.Range("AA137").Value = -1
$ws.Range("AB137").Value = 0.875
$ws.Range("AC137").Value = -0.5
$ws.Range("AD137").Value = 0.5

# Row 139 <- data from former row 140
$ws.Range("B139").Value = 7528849
$ws.Range("C139").Value = "Ecuador LigaPro Serie A"
$ws.Range("D139").Value = 45262.70833333334
$ws.Range("E139").Value = "Guayaquil City"
$ws.Range("F139").Value = "Gualaceo SC"
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 2
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 1
$ws.Range("K139").Value = "A"
$ws.Range("L139").Value = 1.833
$ws.Range("M139").Value = 3.5
$ws.Range("N139").Value = 3.75
$ws.Range("O139").Value = 2.15
$ws.Range("P139").Value = 3.4
$ws.Range("Q139").Value = 3
$ws.Range("R139").Value = -0.25
$ws.Range("S139").Value = 1.825
$ws.Range("T139").Value = 1.975
$ws.Range("U139").Value = 2.5
$ws.Range("V139").Value = 1.85
$ws.Range("W139").Value = 1.95
$ws.Range("X139").Value = -1
$ws.Range("Y139").Value = -1
$ws.Range("Z139").Value = 2
$ws.Range("AA139").Value = -1
$ws.Range("AB139").Value = 0.9750000000000001
$ws.Range("AC139").Value = -1
$ws.Range("AD139").Value = 0.95

# Row 140 <- data from former row 139
$ws.Range("B140").Value = 7528859
$ws.Range("C140").Value = "Ecuador LigaPro Serie A"
$ws.Range("D140").Value = 45262.70833333334
$ws.Range("E140").Value = "Club Atletico Libertad"
$ws.Range("F140").Value = "Cumbaya FC"
$ws.Range("G140").Value = 3
$ws.Range("H140").Value = 1
$ws.Range("I140").Value = 2
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = "H"
$ws.Range("L140").Value = 1.727
$ws.Range("M140").Value = 3.5
$ws.Range("N140").Value = 4.333
$ws.Range("O140").Value = 1.4
$ws.Range("P140").Value = 4.2
$ws.Range("Q140").Value = 7
$ws.Range("R140").Value = -1.25
$ws.Range("S140").Value = 2
$ws.Range("T140").Value = 1.8
$ws.Range("U140").Value = 2.5
$ws.Range("V140").Value = 1.95
$ws.Range("W140").Value = 1.85
$ws.Range("X140").Value = 0.3999999999999999
$ws.Range("Y140").Value = -1
$ws.Range("Z140").Value = -1
$ws.Range("AA140").Value = 1
$ws.Range("AB140").Value = -1
$ws.Range("AC140").Value = 0.95
$ws.Range("AD140").Value = -1

# Row 144 <- data from former row 145
$ws.Range("B144").Value = 7528852
$ws.Range("C144").Value = "Ecuador LigaPro Serie A"
$ws.Range("D144").Value = 45263.83333333334
$ws.Range("E144").Value = "Delfin SC"
$ws.Range("F144").Value = "Tecnico Universitario"
$ws.Range("G144").Value = 2
$ws.Range("H144").Value = 2
$ws.Range("I144").Value = 1
$ws.Range("J144").Value = 0
$ws.Range("K144").Value = "D"
$ws.Range("L144").Value = 2.1
$ws.Range("M144").Value = 3.4
$ws.Range("N144").Value = 3.1
$ws.Range("O144").Value = 2.1
$ws.Range("P144").Value = 3.4
$ws.Range("Q144").Value = 3.1
$ws.Range("R144").Value = -0.25
$ws.Range("S144").Value = 1.8
$ws.Range("T144").Value = 2
$ws.Range("U144").Value = 2.25
$ws.Range("V144").Value = 1.9
$ws.Range("W144").Value = 1.9
$ws.Range("X144").Value = -1
$ws.Range("Y144").Value = 2.4
$ws.Range("Z144").Value = -1
$ws.Range("AA144").Value = -0.5
$ws.Range("AB144").Value = 0.5
$ws.Range("AC144").Value = 0.8999999999999999
$ws.Range("AD144").Value = -1

# Row 145 <- data from former row 144
$ws.Range("B145").Value = 7528857
$ws.Range("C145").Value = "Ecuador LigaPro Serie A"
$ws.Range("D145").Value = 45263.83333333334
$ws.Range("E145").Value = "Universidad Catolica del Ecuador"
$ws.Range("F145").Value = "Barcelona Guayaquil"
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 1
$ws.Range("I145").Value = 0
$ws.Range("J145").Value = 0
$ws.Range("K145").Value = "A"
$ws.Range("L145").Value = 1.533
$ws.Range("M145").Value = 4
$ws.Range("N145").Value = 5.5
$ws.Range("O145").Value = 1.5
$ws.Range("P145").Value = 4.333
$ws.Range("Q145").Value = 5.25
$ws.Range("R145").Value = -1
$ws.Range("S145").Value = 1.8
$ws.Range("T145").Value = 2
$ws.Range("U145").Value = 3
$ws.Range("V145").Value = 1.975
$ws.Range("W145").Value = 1.825
$ws.Range("X145").Value = -1
$ws.Range("Y145").Value = -1
$ws.Range("Z145").Value = 4.25
$ws.Range("AA145").Value = -1
$ws.Range("AB145").Value = 1
$ws.Range("AC145").Value = -1
$ws.Range("AD145").Value = 0.825
